$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 203, pushing the existing rows 203-212 down to 204-213.
$ws.Rows.Item(203).Insert()

# Populate the newly inserted row 203 with the new weekly price observation.
$ws.Cells.Item(203, 1).Value = 9
$ws.Cells.Item(203, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(203, 3).Value = "Metropolitana"
$ws.Cells.Item(203, 4).Value = 44706
$ws.Cells.Item(203, 5).Value = 13
$ws.Cells.Item(203, 6).Value = 100112026
$ws.Cells.Item(203, 7).Value = "Haba"
$ws.Cells.Item(203, 8).Value = "Sin especificar"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 45
$ws.Cells.Item(203, 11).Value = 21000
$ws.Cells.Item(203, 12).Value = 22000
$ws.Cells.Item(203, 13).Value = 21556
$ws.Cells.Item(203, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(203, 15).Value = "Región Metropolitana"
$ws.Cells.Item(203, 16).Value = 862
$ws.Cells.Item(203, 17).Value = 25
$ws.Cells.Item(203, 18).Value = "Hortaliza"
